# Command Strings List - update to match development progress
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section header text tweaks (now show the SubSystem code ranges) ---
$ws.Range("A4").Value = "Rover Systems (8000-9000)"
$ws.Range("A9").Value = "Drive Train (9000-10000)"

# --- Existing "Receive GPS Stream" row: port moved from 9001 to 8001, and it now
#     also documents an example response string ---
$ws.Range("C5").Value = 8001
$ws.Range("J5").Value = "0.000000000,0.0000000000,0.00,0.00,0.0,0"

# --- New row: "Receive Temperature Stream" command, inserted directly under the
#     GPS stream row, mirroring its layout/format ---
$ws.Range("B6").Value = "Receive Temperature Stream"
$ws.Range("C6").Value = 8002
$ws.Range("D6").Value = "N/A"
$ws.Range("E6").Value = "N/A"
$ws.Range("F6").Value = "N/A"
$ws.Range("G6").Value = "Simply open the conenction to receive data"
$ws.Range("I6").Value = "°C"
$ws.Range("J6").Value = 48.234

# --- Column widths grew to fit the new, longer text ---
$ws.Columns.Item(1).ColumnWidth = 32.42
$ws.Columns.Item(2).ColumnWidth = 26.42
$ws.Columns.Item(10).ColumnWidth = 39.34

# --- Selection moved ---
$ws.Range("B12").Select() | Out-Null
